# Applies the two changes described by the commit:
#   1. The cached "datetimeFigureOut" field text on the slide master and
#      every slide layout is restamped from 2021/12/27 -> 2022/6/4
#      (this happens automatically in real PowerPoint whenever the file is
#      re-saved on a different day; here we do it explicitly).
#   2. On slide 1, the mis-spelled run split "l" + "ogisitic" inside the
#      "TextBox 2" shape is corrected/merged into a single "logistic" run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached date field wherever it still shows the old date.
# ---------------------------------------------------------------------
$oldDate = "2021/12/27"
$newDate = "2022/6/4"

function Update-DatePlaceholders {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.HasText) {
                if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
                    $shape.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Fix "l" + "ogisitic" -> "logistic" in TextBox 2 on slide 1.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.Name -eq "TextBox 2") {
        $tr = $shape.TextFrame.TextRange
        # "l" (1 char) + "ogisitic" (8 chars) = first 9 characters of the
        # text range; replacing them in one go merges the two runs into a
        # single run carrying the first run's formatting.
        $fixRange = $tr.Characters(1, 9)
        if ($fixRange.Text -eq "logisitic") {
            $fixRange.Text = "logistic"
        }
    }
}
